# Báo cáo cá nhân - NV-11 Đỗ Thị Huyền Trân 8-2024
# Adds "Đơn sale phụ" sheet between "Đơn sale chính" and "Lương",
# fills both order sheets with their August-2024 rows,
# and updates the "Lương" report figures for CẦN THƠ.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert the new "Đơn sale phụ" worksheet right after
#    "Đơn sale chính" (so tab order becomes chính / phụ / Lương).
# ---------------------------------------------------------------
$wsChinh = $wb.Worksheets.Item(1)
$wsPhu = $wb.Worksheets.Add([System.Type]::Missing, $wsChinh)
$wsPhu.Name = "Đơn sale phụ"

$wsChinh = $wb.Worksheets.Item(1)
$wsPhu = $wb.Worksheets.Item(2)
$wsLuong = $wb.Worksheets.Item(3)

# Match the page-setup/outline defaults used by the sibling sheets.
$wsPhu.Outline.SummaryRow = 1
$wsPhu.Outline.SummaryColumn = 1
$wsPhu.PageSetup.LeftMargin = 54
$wsPhu.PageSetup.RightMargin = 54
$wsPhu.PageSetup.TopMargin = 72
$wsPhu.PageSetup.BottomMargin = 72
$wsPhu.PageSetup.HeaderMargin = 36
$wsPhu.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------
# Helper data: common column headers (cols A..N)
# ---------------------------------------------------------------
$colCount = 14

# ---------------------------------------------------------------
# 2. "Đơn sale chính" sheet content
# ---------------------------------------------------------------
$headersChinh = @("Tiền tố","Mã dịch vụ","Ngày thực hiện","Cơ sở","Khách hàng","Nguồn khách","Tên dịch vụ","Đơn giá gốc","Sale phụ","Upsale","Đơn giá","Đã thanh toán","Tỉ lệ chiết khấu sale chính","Chiết khấu sale chính")
for ($i = 0; $i -lt $headersChinh.Length; $i++) {
    $wsChinh.Cells.Item(1, $i + 1).Value = $headersChinh[$i]
}

$wsChinh.Cells.Item(2, 1).Value = "HD-LUXURY"
$wsChinh.Cells.Item(2, 2).Value = 615
$wsChinh.Cells.Item(2, 3).NumberFormat = "@"
$wsChinh.Cells.Item(2, 3).Value = "08-01-2024"
$wsChinh.Cells.Item(2, 4).Value = "CẦN THƠ"
$wsChinh.Cells.Item(2, 5).Value = "Nguyễn Thị Mỹ Duyên"
$wsChinh.Cells.Item(2, 6).Value = "Khách cũ"
$wsChinh.Cells.Item(2, 7).Value = "Tiêm Filler"
$wsChinh.Cells.Item(2, 8).Value = 2100000
$wsChinh.Cells.Item(2, 9).Value = 0
$wsChinh.Cells.Item(2, 10).Value = 0
$wsChinh.Cells.Item(2, 11).Value = 2100000
$wsChinh.Cells.Item(2, 12).Value = 2100000
$wsChinh.Cells.Item(2, 13).Value = 0.1
$wsChinh.Cells.Item(2, 14).Value = 210000

$wsChinh.Cells.Item(3, 1).Value = "Tổng"
$wsChinh.Cells.Item(3, 2).Value = 1
$wsChinh.Cells.Item(3, 8).Value = 2100000
$wsChinh.Cells.Item(3, 10).Value = 0
$wsChinh.Cells.Item(3, 11).Value = 2100000
$wsChinh.Cells.Item(3, 12).Value = 2100000
$wsChinh.Cells.Item(3, 13).Value = 0
$wsChinh.Cells.Item(3, 14).Value = 210000

# ---------------------------------------------------------------
# 3. "Đơn sale phụ" sheet content
# ---------------------------------------------------------------
$headersPhu = @("Tiền tố","Mã dịch vụ","Ngày thực hiện","Cơ sở","Khách hàng","Nguồn khách","Tên dịch vụ","Đơn giá gốc","Sale phụ","Upsale","Đơn giá","Đã thanh toán","Tỉ lệ chiết khấu sale phụ","Chiết khấu sale phụ")
for ($i = 0; $i -lt $headersPhu.Length; $i++) {
    $wsPhu.Cells.Item(1, $i + 1).Value = $headersPhu[$i]
}

$wsPhu.Cells.Item(2, 1).Value = "HD-LUXURY"
$wsPhu.Cells.Item(2, 2).Value = 614
$wsPhu.Cells.Item(2, 3).NumberFormat = "@"
$wsPhu.Cells.Item(2, 3).Value = "08-01-2024"
$wsPhu.Cells.Item(2, 4).Value = "CẦN THƠ"
$wsPhu.Cells.Item(2, 5).Value = "Trần Nguyễn Yến Linh"
$wsPhu.Cells.Item(2, 6).Value = "Khách cũ"
$wsPhu.Cells.Item(2, 7).Value = "Cắt mí"
$wsPhu.Cells.Item(2, 8).Value = 0
$wsPhu.Cells.Item(2, 9).Value = "Đỗ Thị Huyền Trân"
$wsPhu.Cells.Item(2, 10).Value = 6000000
$wsPhu.Cells.Item(2, 11).Value = 6000000
$wsPhu.Cells.Item(2, 12).Value = 6000000
$wsPhu.Cells.Item(2, 13).Value = 0.04
$wsPhu.Cells.Item(2, 14).Value = 240000

$wsPhu.Cells.Item(3, 1).Value = "Tổng"
$wsPhu.Cells.Item(3, 2).Value = 1
$wsPhu.Cells.Item(3, 8).Value = 0
$wsPhu.Cells.Item(3, 10).Value = 6000000
$wsPhu.Cells.Item(3, 11).Value = 6000000
$wsPhu.Cells.Item(3, 12).Value = 6000000
$wsPhu.Cells.Item(3, 13).Value = 0
$wsPhu.Cells.Item(3, 14).Value = 240000

# ---------------------------------------------------------------
# 4. "Lương" sheet: update the CẦN THƠ figures now that orders exist
# ---------------------------------------------------------------
$wsLuong.Cells.Item(2, 2).Value = 2
$wsLuong.Cells.Item(3, 2).Value = 70000
$wsLuong.Cells.Item(4, 2).Value = 357142.8571428572
$wsLuong.Cells.Item(5, 2).Value = 210000
$wsLuong.Cells.Item(6, 2).Value = 240000
$wsLuong.Cells.Item(32, 2).Value = 877142.8571428572
$wsLuong.Cells.Item(35, 1).Value = "Tổng lương tại HỆ THỐNG"
$wsLuong.Cells.Item(35, 2).Value = 877142.8571428572
